$d = $word.ActiveDocument

# 1. Remove ", while possibly being a bit less accurate than some other models,"
$d.Content.Find.Execute(", while possibly being a bit less accurate than some other models,", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2. Insert " (at least in comparison to the previous homework)" after "train our models" (before the period)
$d.Content.Find.Execute("train our models. In my", $true, $false, $false, $false, $false, $true, 1, $false, "train our models (at least in comparison to the previous homework). In my", 2)

# 3. Insert "every iteration " before "until the predictions"
$d.Content.Find.Execute("update weights until the predictions", $true, $false, $false, $false, $false, $true, 1, $false, "update weights every iteration until the predictions", 2)

# 4. Insert " (aggressive)" before ". This leads"
$d.Content.Find.Execute("precisely correct. This leads", $true, $false, $false, $false, $false, $true, 1, $false, "precisely correct (aggressive). This leads", 2)
